# Auto-generated: updates cryptos list cell values to match target snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $ref, $val) {
    # Excel auto-converts strings that look like plain numbers into Number
    # cells (dropping things like trailing zeros). The source data stores
    # these as literal text, so force a text number-format first whenever the
    # new value would otherwise be re-interpreted as a number.
    if ($val -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        $ws.Range($ref).NumberFormat = "@"
    }
    $ws.Range($ref).Value = $val
}

Set-CellText $ws "D2" "69.334.26"
Set-CellText $ws "E2" "  -0.20%  "
Set-CellText $ws "D3" "3.686.60"
Set-CellText $ws "E3" "  -0.17%  "
Set-CellText $ws "D5" "677.82"
Set-CellText $ws "E5" "  -1.63%  "
Set-CellText $ws "D6" "159.08"
Set-CellText $ws "E6" "  -1.80%  "
Set-CellText $ws "E7" "  -0.03%  "
Set-CellText $ws "E8" "  -0.44%  "
Set-CellText $ws "E9" "  -1.09%  "
Set-CellText $ws "D10" "7.16"
Set-CellText $ws "E10" "  -2.82%  "
Set-CellText $ws "D11" "0.442"
Set-CellText $ws "E11" "  +0.57%  "
Set-CellText $ws "E12" "  -2.26%  "
Set-CellText $ws "D13" "4.308.29"
Set-CellText $ws "E13" "  -0.15%  "
Set-CellText $ws "D14" "32.37"
Set-CellText $ws "E14" "  -2.11%  "
Set-CellText $ws "D15" "3.707.08"
Set-CellText $ws "E15" "  +0.46%  "
Set-CellText $ws "D16" "69.316.78"
Set-CellText $ws "E16" "  -0.21%  "
Set-CellText $ws "E17" "  +2.82%  "
Set-CellText $ws "E18" "  -0.46%  "
Set-CellText $ws "D19" "6.49"
Set-CellText $ws "E19" "  -0.61%  "
Set-CellText $ws "D20" "468.22"
Set-CellText $ws "E20" "  -1.78%  "
Set-CellText $ws "E21" "  -0.91%  "
Set-CellText $ws "E22" "  -0.76%  "
Set-CellText $ws "D23" "79.99"
Set-CellText $ws "E23" "  +0.01%  "
Set-CellText $ws "D24" "3.833.24"
Set-CellText $ws "E24" "  -0.10%  "
Set-CellText $ws "E25" "  -0.08%  "
Set-CellText $ws "E26" "  -4.98%  "
Set-CellText $ws "D27" "10.91"
Set-CellText $ws "E27" "  -2.87%  "
Set-CellText $ws "E28" "  -1.36%  "
Set-CellText $ws "E29" "  -0.78%  "
Set-CellText $ws "E30" "  -3.53%  "
Set-CellText $ws "E31" "  -3.14%  "
Set-CellText $ws "E33" "  -3.00%  "
Set-CellText $ws "D34" "26.91"
Set-CellText $ws "E34" "  +0.11%  "
Set-CellText $ws "D35" "3.676.44"
Set-CellText $ws "E35" "  +0.56%  "
Set-CellText $ws "D36" "0.160"
Set-CellText $ws "E36" "  -3.93%  "
Set-CellText $ws "D37" "8.31"
Set-CellText $ws "E37" "  -0.42%  "
Set-CellText $ws "D38" "6.30"
Set-CellText $ws "E38" "  +0.17%  "
Set-CellText $ws "E39" "  -0.01%  "
Set-CellText $ws "B40" "FirstDigitalUSD"
Set-CellText $ws "C40" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText $ws "D40" "1.00"
Set-CellText $ws "E40" "  -0.04%  "
Set-CellText $ws "B41" "Stacks"
Set-CellText $ws "C41" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-CellText $ws "D41" "2.24"
Set-CellText $ws "E41" "  -3.15%  "
Set-CellText $ws "E42" "  -1.60%  "
Set-CellText $ws "D43" "170.62"
Set-CellText $ws "E43" "  +4.34%  "
Set-CellText $ws "D44" "0.942"
Set-CellText $ws "E44" "  -1.08%  "
Set-CellText $ws "D45" "47.28"
Set-CellText $ws "E45" "  -2.18%  "
Set-CellText $ws "D46" "28.15"
Set-CellText $ws "E46" "  -7.10%  "
Set-CellText $ws "E47" "  -0.85%  "
Set-CellText $ws "E48" "  -2.50%  "
Set-CellText $ws "E49" "  -2.16%  "
Set-CellText $ws "D50" "1.29"
Set-CellText $ws "E50" "  -2.25%  "
Set-CellText $ws "E51" "  -2.78%  "
